$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 68628250
$ws.Range("I43").Value = 250000260
$ws.Range("J43").Value = 12821480
$ws.Range("K43").Value = 250000260
$ws.Range("L43").Value = 12821480
$ws.Range("M43").Value = -250000191
$ws.Range("N43").Value = -12821618

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1372.5
$ws.Range("J48").Value = 1372.5
$ws.Range("L48").Value = 4117.5
$ws.Range("N48").Value = -4701.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 1372.5
$ws.Range("J56").Value = 1372.5
$ws.Range("L56").Value = 4117.5
$ws.Range("N56").Value = -5185.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3980.1353
$ws.Range("I69").Value = 3447.647
$ws.Range("J69").Value = 10015
$ws.Range("K69").Value = 10342.941
$ws.Range("L69").Value = 30045
$ws.Range("M69").Value = -9468.940999999999
$ws.Range("N69").Value = -31793

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3980.1353
$ws.Range("I72").Value = 3447.647
$ws.Range("J72").Value = 10015
$ws.Range("K72").Value = 31028.823
$ws.Range("L72").Value = 90135
$ws.Range("M72").Value = -26660.823
$ws.Range("N72").Value = -98871

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 629.4286
$ws.Range("I135").Value = 528.9787
$ws.Range("J135").Value = 2990
$ws.Range("K135").Value = 4760.8083
$ws.Range("L135").Value = 26910
$ws.Range("M135").Value = -2225.8083
$ws.Range("N135").Value = -31980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1952.2069
$ws.Range("I137").Value = 1610.6316
$ws.Range("J137").Value = 2601.2
$ws.Range("K137").Value = 4831.8948
$ws.Range("L137").Value = 7803.599999999999
$ws.Range("M137").Value = -2281.8948
$ws.Range("N137").Value = -12903.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1873.421
$ws.Range("I138").Value = 1526
$ws.Range("K138").Value = 4578
$ws.Range("M138").Value = 562

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2338.9722
$ws.Range("I141").Value = 1403.25
$ws.Range("J141").Value = 3508.625
$ws.Range("K141").Value = 4209.75
$ws.Range("L141").Value = 10525.875
$ws.Range("M141").Value = 970.25
$ws.Range("N141").Value = -20885.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1020.6087
$ws.Range("I61").Value = 642.8627300000001
$ws.Range("J61").Value = 2090.889
$ws.Range("K61").Value = 642.8627300000001
$ws.Range("L61").Value = 2090.889
$ws.Range("M61").Value = -430.8627300000001
$ws.Range("N61").Value = -2514.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 943.7778
$ws.Range("I74").Value = 843
$ws.Range("K74").Value = 843
$ws.Range("M74").Value = 31

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 943.7778
$ws.Range("I77").Value = 843
$ws.Range("K77").Value = 4215
$ws.Range("M77").Value = 153

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1020.6087
$ws.Range("I136").Value = 642.8627300000001
$ws.Range("J136").Value = 2090.889
$ws.Range("K136").Value = 1928.58819
$ws.Range("L136").Value = 6272.667
$ws.Range("M136").Value = 621.4118099999998
$ws.Range("N136").Value = -11372.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2885.3225
$ws.Range("I20").Value = 3023.0527
$ws.Range("K20").Value = 3023.0527
$ws.Range("M20").Value = -2776.0527

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22286.307
$ws.Range("I134").Value = 28655
$ws.Range("J134").Value = 2649.5
$ws.Range("K134").Value = 85965
$ws.Range("L134").Value = 7948.5
$ws.Range("M134").Value = -83430
$ws.Range("N134").Value = -13018.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4168671.8
$ws.Range("I31").Value = 1295.625
$ws.Range("J31").Value = 20838176
$ws.Range("K31").Value = 1295.625
$ws.Range("L31").Value = 20838176
$ws.Range("M31").Value = -1000.625
$ws.Range("N31").Value = -20838766

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4168671.8
$ws.Range("I34").Value = 1295.625
$ws.Range("J34").Value = 20838176
$ws.Range("K34").Value = 1295.625
$ws.Range("L34").Value = 20838176
$ws.Range("M34").Value = -1093.625
$ws.Range("N34").Value = -20838580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1329.2903
$ws.Range("I134").Value = 1315.6923
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 3947.0769
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -1412.0769
$ws.Range("N134").Value = -9270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1257.9
$ws.Range("I5").Value = 504
$ws.Range("J5").Value = 1446.375
$ws.Range("K5").Value = 1512
$ws.Range("L5").Value = 4339.125
$ws.Range("M5").Value = -1400
$ws.Range("N5").Value = -4563.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1212.2632
$ws.Range("J122").Value = 966.2727
$ws.Range("L122").Value = 8696.454299999999
$ws.Range("N122").Value = -13596.4543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1025027.4
$ws.Range("I131").Value = 5977
$ws.Range("J131").Value = 3572653.5
$ws.Range("K131").Value = 17931
$ws.Range("L131").Value = 10717960.5
$ws.Range("M131").Value = -12891
$ws.Range("N131").Value = -10728040.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1257.9
$ws.Range("I135").Value = 504
$ws.Range("J135").Value = 1446.375
$ws.Range("K135").Value = 4536
$ws.Range("L135").Value = 13017.375
$ws.Range("M135").Value = -2001
$ws.Range("N135").Value = -18087.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 48.15
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 62.25
$ws.Range("K2").Value = 27
$ws.Range("L2").Value = 62.25
$ws.Range("M2").Value = 86
$ws.Range("N2").Value = -288.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3511111
$ws.Range("I7").Value = 3825000
$ws.Range("K7").Value = 3825000
$ws.Range("M7").Value = -3824888

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 3511111
$ws.Range("I8").Value = 3825000
$ws.Range("K8").Value = 3825000
$ws.Range("M8").Value = -3824861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11864560
$ws.Range("I70").Value = 13787891
$ws.Range("J70").Value = 4016.6667
$ws.Range("K70").Value = 13787891
$ws.Range("L70").Value = 4016.6667
$ws.Range("M70").Value = -13787621
$ws.Range("N70").Value = -4556.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11864560
$ws.Range("I73").Value = 13787891
$ws.Range("J73").Value = 4016.6667
$ws.Range("K73").Value = 13787891
$ws.Range("L73").Value = 4016.6667
$ws.Range("M73").Value = -13786955
$ws.Range("N73").Value = -5888.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4064
$ws.Range("I80").Value = 2850.8
$ws.Range("K80").Value = 2850.8
$ws.Range("M80").Value = -1852.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4064
$ws.Range("I83").Value = 2850.8
$ws.Range("K83").Value = 14254
$ws.Range("M83").Value = -9262

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 943.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 35717024
$ws.Range("I122").Value = 76926230
$ws.Range("K122").Value = 230778690
$ws.Range("M122").Value = -230776240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2879.111
$ws.Range("I126").Value = 3244.5715
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 9733.7145
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -7263.7145
$ws.Range("N126").Value = -9740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 17000
$ws.Range("J50").Value = 17000
$ws.Range("L50").Value = 17000
$ws.Range("N50").Value = -18274

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1014.60376
$ws.Range("I132").Value = 971
$ws.Range("K132").Value = 2913
$ws.Range("M132").Value = -383

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 91299.664
$ws.Range("J135").Value = 91299.664
$ws.Range("L135").Value = 91299.664
$ws.Range("N135").Value = -101439.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3186.9812
$ws.Range("I136").Value = 3483.2954
$ws.Range("K136").Value = 10449.8862
$ws.Range("M136").Value = -7899.886200000001
